# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# The simulator recomputed the per-game strikeout totals (K) for this
# player's 2021 log; write the refreshed values back into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 3
    11 = 2
    12 = 1
    13 = 1
    14 = 3
    15 = 1
    16 = 0
    17 = 2
    18 = 3
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 0
    24 = 2
    26 = 0
    27 = 1
    28 = 3
    29 = 1
    30 = 1
    31 = 1
    32 = 2
    33 = 0
    34 = 0
    35 = 0
    36 = 1
    37 = 1
    38 = 1
    40 = 2
    42 = 1
    43 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
